# Rename the two worksheets ("sheet names in xlsx")
$wb = $excel.ActiveWorkbook

$wsLipids = $wb.Worksheets.Item("Sheet3")
$wsLipids.Name = "lipids"

$wsFiles = $wb.Worksheets.Item("Sheet4")
$wsFiles.Name = "files"

# The "lipids" sheet had a helper column J with a "Theoretical mass" header
# and per-row computed values that are no longer needed - clear them out
# (this also drops the now-unused "Theoretical mass" shared string).
$wsLipids.Range("J1:J39").ClearContents()

# Restore a plain cell selection on the lipids sheet (was J2:J41).
$wsLipids.Activate()
$wsLipids.Range("E11").Select()
